$wb = $excel.ActiveWorkbook

# Rename the empty "Sheet1" worksheet to "createpotential" and use it to
# build out the new "Create Potential" test-data sheet.
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "createpotential"

# Header row (row 1) - plain text values, General-formatted cells.
$ws.Range("A1").Value = "Login tilte"
$ws.Range("B1").Value = "Hometitle"
$ws.Range("C1").Value = "Potential page"
$ws.Range("D1").Value = "NewPotential"
$ws.Range("E1").Value = "Potential name"
$ws.Range("F1").Value = "ClosingDate"

# Data row (row 2) - stored as text (numFmtId 49 "@") to match existing
# strings used elsewhere in the workbook.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "Zoho CRM - Sign in"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "Zoho CRM - Home Page"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "Zoho CRM - Displaying Custom View Details"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "Zoho CRM - Create Potential"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "jack"

# Closing date is entered as a literal text string (apostrophe-prefixed)
# rather than a real date, matching the source workbook's "quote prefix"
# text cell.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "'11/30/2019"

# Make this newly populated sheet the active tab / selection, mirroring
# the state Excel saved the workbook in.
$ws.Activate()
$ws.Range("A1:F2").Select()
